# Apply the "Lot's of clean up. New featurs added to HTML template." edit
# to the Website sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Website")

# The "new customer" blurb is rewritten and moved to be the last shared
# string (row 6 / C6): "Cards that create new customers!" becomes
# "Your cards can create new customers!"
$ws.Range("C6").Value = "Your cards can create new customers!"

# The QR-code row numbering is renumbered from 9/10 to 5/6.
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 6

# Update the sheet's remembered selection to B3:B8 (active cell B3).
$ws.Range("B3:B8").Select()
